$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Semestre ideal" value from "EP-8" to "EA-8,EP-8" (columns B and C, row 9)
$ws.Range("B9").Value = "EA-8,EP-8"
$ws.Range("C9").Value = "EA-8,EP-8"

# Remove the "Requisitos:" rows (24 and 25) entirely
$ws.Rows("24:25").Delete()
